# Insert a new data row at row 231 (shifts existing rows 231-284 down to 232-285)
# and populate it with a new price observation, matching the target diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(231).Insert()

$ws.Cells.Item(231, 1).Value = 7
$ws.Cells.Item(231, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(231, 3).Value = "Ñuble"
$ws.Cells.Item(231, 4).Value = 44932
$ws.Cells.Item(231, 5).Value = 16
$ws.Cells.Item(231, 6).Value = 100112043
$ws.Cells.Item(231, 7).Value = "Pepino ensalada"
$ws.Cells.Item(231, 8).Value = "Sin especificar"
$ws.Cells.Item(231, 9).Value = "Primera"
$ws.Cells.Item(231, 10).Value = 120
$ws.Cells.Item(231, 11).Value = 11000
$ws.Cells.Item(231, 12).Value = 12000
$ws.Cells.Item(231, 13).Value = 11500
$ws.Cells.Item(231, 14).Value = "`$/caja 80 unidades"
$ws.Cells.Item(231, 15).Value = "Región del Maule"
$ws.Cells.Item(231, 16).Value = 144
$ws.Cells.Item(231, 17).Value = 80
$ws.Cells.Item(231, 18).Value = "Hortaliza"
